# Template_for_Adding_Courses_to_Sustainability_Course_Finder.xlsx update
# - Split the single "course_title" column into a generic title plus a new
#   "section_name" column (inserted as column F, pushing instructor..course_level
#   from F:L to G:M).
# - Rename course_desc -> course_description and drop the redundant
#   "Special Topics: " prefix from the long descriptions.
# - Refresh instructor list, section numbers, semester codes, and add a new
#   trailing "session" column with sample section counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Preserve formatting that must travel with specific values as the
#    instructor..course_level block shifts one column to the right.
#    (Do this before the value-shift, while the "donor" cells still carry
#    their original look.)
# ---------------------------------------------------------------------
# The old "date" formatting on G2 belongs with the date once it moves to H2.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# The old "highlighted" course_title formatting on E2 belongs with the long
# title text once it becomes the section_name value in F2.
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Shift the instructor..course_level values one column to the right
#    (F:L -> G:M), freeing up column F for the new section_name column.
# ---------------------------------------------------------------------
$ws.Range("G1:M4").Value = $ws.Range("F1:L4").Value()
$ws.Range("F1:F4").ClearContents()

# G2 now holds the instructor text (it used to hold the date), so restore
# the plain row-2 formatting there (matching its sibling D2).
$ws.Range("D2").Copy()
$ws.Range("G2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# New trailing M2:M4 cells need the base formatting used elsewhere in their
# row (copied from the already-shifted L column).
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "section_name"
$ws.Range("J1").Value = "course_description"
$ws.Range("N1").Value = "session"

# ---------------------------------------------------------------------
# 4) Row 2 - ENGR-499 example row
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Update Description"
$ws.Range("C2").Value = "ENGR"
$ws.Range("E2").Value = "Special Topics"
$ws.Range("F2").Value = "Systems Thinking for Sustainability"
$ws.Range("G2").Value = "Albright, Julie;Maby, Edward, W"
$ws.Range("J2").Value = "This course will grow your understanding about the foundational systems shaping our world - energy, communications, transportation, water, waste, mining, and others - emphasizing their intersection (a ""systems approach"") with the social, environmental, business, policy, and technical arenas."

# ---------------------------------------------------------------------
# 5) Row 3 - AME-599 Sustainable Aerospace row
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Update Description"
$ws.Range("E3").Value = "Special Topics"
$ws.Range("F3").Value = "Sustainable Aerospace"
$ws.Range("H3").Value = 29085
$ws.Range("I3").Value = "F24"
$ws.Range("J3").Value = "This course presents the history and current developments in the field of sustainable aerospace, covering both aviation and space topics. Topics will include alternative fuels (biofuels, synthetic fuels, methane, alcohols, and hydrogen), electric and hybrid electric aircraft, lifecycle environmental impact and analysis, space debris, the environmental impact of rocket launches and space debris reentry, and aerospace technology based ideas to mitigate climate change.  The advantages and challenges of each type of potential sustainable aerospace technology will be discussed, evaluated, and compared. "
$ws.Range("N3").Value = 48

# ---------------------------------------------------------------------
# 6) Row 4 - AME-499 Sustainable Aviation row
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Add"
$ws.Range("E4").Value = "Special Topics"
$ws.Range("F4").Value = "Sustainable Aviation"
$ws.Range("H4").Value = 29035
$ws.Range("I4").Value = "F24"
$ws.Range("J4").Value = "This course presents the history and current developments in the field of sustainable aviation. Topics will include alternative fuels (biofuels, synthetic fuels, methane, alcohols, and hydrogen), electric and hybrid electric aircraft, and lifecycle environmental impact and analysis.  The advantages and challenges of each type of potential sustainable aviation technology will be discussed, evaluated, and compared.  The target audience for this course is undergraduate and graduate students with an interest or background in aerospace engineering or sustainable energy who want to learn about sustainable aviation."
$ws.Range("N4").Value = 966

# Selection/view tweak to match the saved workbook
$ws.Range("A4").Select()
